$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name (tab) to reflect new "through" date
$ws.Name = "Through 2022-04-29"

# Update the April row label
$ws.Range("A5").Value = "April (through 04-29)"

# Update April row (row 5) values
$ws.Range("B5").Value = 21
$ws.Range("C5").Value = 32
$ws.Range("D5").Value = 61
$ws.Range("E5").Value = 47
$ws.Range("F5").Value = 43
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 98
$ws.Range("I5").Value = 116

# Update Total row (row 6) values
$ws.Range("B6").Value = 87
$ws.Range("C6").Value = 160
$ws.Range("D6").Value = 250
$ws.Range("E6").Value = 244
$ws.Range("F6").Value = 153
$ws.Range("G6").Value = 258
$ws.Range("H6").Value = 521
$ws.Range("I6").Value = 551

$wb.Save()
